# Update the cryptos worksheet with refreshed price / volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2'  = '42.807.14'
    'E2'  = '  +0.95%  '
    'D3'  = '2.294.20'
    'E3'  = '  -0.35%  '
    'E4'  = '  +0.15%  '
    'D5'  = '316.28'
    'E5'  = '  -0.24%  '
    'D6'  = '104.15'
    'E6'  = '  -0.09%  '
    'D7'  = '0.624'
    'E7'  = '  -1.13%  '
    'E8'  = '  +0.14%  '
    'E9'  = '  -1.83%  '
    'D10' = '39.42'
    'E10' = '  -1.95%  '
    'E11' = '  -0.59%  '
    'D12' = '8.48'
    'E12' = '  +1.58%  '
    'E13' = '  +2.09%  '
    'E14' = '  +4.18%  '
    'D15' = '15.30'
    'E15' = '  -0.08%  '
    'D16' = '2.642.75'
    'E16' = '  -0.25%  '
    'D17' = '2.309.75'
    'E17' = '  +0.25%  '
    'D18' = '42.732.90'
    'E18' = '  +0.53%  '
    'E19' = '  -0.79%  '
    'D20' = '14.04'
    'E20' = '  +25.75%  '
    'E21' = '  -0.59%  '
    'E22' = '  +0.99%  '
    'E23' = '  -0.44%  '
    'D24' = '263.94'
    'E24' = '  -4.47%  '
    'E25' = '  -3.17%  '
    'E26' = '  +0.28%  '
    'E27' = '  +0.45%  '
    'D28' = '7.12'
    'E28' = '  +20.95%  '
    'E29' = '  +0.06%  '
    'D30' = '22.36'
    'E30' = '  -1.86%  '
    'D31' = '37.60'
    'E31' = '  +4.85%  '
    'D32' = '166.68'
    'E32' = '  +1.09%  '
    'E33' = '  -0.54%  '
    'E34' = '  -4.47%  '
    'E35' = '  -0.62%  '
    'E36' = '  -1.51%  '
    'E37' = '  -0.72%  '
    'D38' = '0.0349'
    'E38' = '  -6.31%  '
    'D39' = '3.82'
    'E39' = '  +1.28%  '
    'E40' = '  -2.80%  '
    'E41' = '  +3.94%  '
    'E42' = '  +1.34%  '
    'D43' = '69.28'
    'E43' = '  -1.28%  '
    'E44' = '  -0.06%  '
    'D45' = '91.93'
    'E45' = '  -3.39%  '
    'D46' = '12.12'
    'E46' = '  +0.40%  '
    'D47' = '114.22'
    'E47' = '  +1.07%  '
    'B48' = 'ordi'
    'C48' = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
    'D48' = '80.56'
    'E48' = '  -2.59%  '
    'B49' = 'Maker'
    'C49' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D49' = '1.724.17'
    'E49' = '  +8.32%  '
    'D50' = '8.78'
    'E50' = '  -1.26%  '
    'D51' = '5.13'
    'E51' = '  +1.38%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force the value to be stored as text (matching the source data, which
    # is all inline/shared-string text even when it looks numeric), then
    # reset the style back to Normal so no formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
